$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped from 2023-09-11
# (45180) to 2023-09-12 (45181) for every data row (rows 2-347).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 347 }

$ws.Range("C2:C$lastRow").Value = 45181
